$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the header row (before the old row 2),
# shifting the existing data rows down.
$ws.Rows("2:4").Insert()

# The insert picks up formatting from the row above (the bold header row);
# reset the new rows back to the default (unstyled) format used by the
# rest of the data rows.
$ws.Rows("2:4").ClearFormats()

# Fill the newly inserted rows (2-4) with their new data.
$newTop = @(
    @(0.08857546001672741, -0.0678060427308082, -0.0534507073462009),
    @(0.0888808965682983, 0.1319468915462494, -0.0287106670439243),
    @(0.3381139039993286, 0.00534507073462, 0.1145372316241264)
)

for ($i = 0; $i -lt $newTop.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTop[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTop[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTop[$i][2]
}

# Append 6 new rows (25-31) with new data after the shifted data (now ending at row 24).
$newBottom = @(
    @(0.1823432743549347, 0.1472185254096985, -0.3480404615402221),
    @(0.5630650520324707, 0.6291912198066711, -0.8999572396278381),
    @(-0.0210748501121997, -0.0826195254921913, -0.0035124751739203),
    @(-0.0372627787292003, 0.0313068442046642, -0.06322455406188961),
    @(0.0158824957907199, 0.0445931628346443, 0.0317649915814399),
    @(0.0148134818300604, -0.0259617734700441, 0.0120645882561802),
    @(0.0215329993516206, 0.0319177098572254, 0.009315694682300001)
)

for ($i = 0; $i -lt $newBottom.Length; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 1).Value = $newBottom[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottom[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottom[$i][2]
}
